$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45859.01041666666, 1),
    @(45859.02083333334, 1),
    @(45859.03125, 1),
    @(45859.04166666666, 1),
    @(45859.05208333334, 1),
    @(45859.0625, 1),
    @(45859.07291666666, 1),
    @(45859.08333333334, 1),
    @(45859.09375, 1),
    @(45859.10416666666, 1),
    @(45859.11458333334, 1),
    @(45859.125, 1),
    @(45859.13541666666, 1),
    @(45859.14583333334, 1),
    @(45859.15625, 1),
    @(45859.16666666666, 1),
    @(45859.17708333334, 1),
    @(45859.1875, 1),
    @(45859.19791666666, 1),
    @(45859.20833333334, 2),
    @(45859.21875, 47),
    @(45859.22916666666, 53),
    @(45859.23958333334, 62),
    @(45859.25, 77),
    @(45859.26041666666, 299),
    @(45859.27083333334, 319),
    @(45859.28125, 352),
    @(45859.29166666666, 394),
    @(45859.30208333334, 796),
    @(45859.3125, 838),
    @(45859.32291666666, 902),
    @(45859.33333333334, 966),
    @(45859.34375, 1397),
    @(45859.35416666666, 1434),
    @(45859.36458333334, 1479),
    @(45859.375, 1530),
    @(45859.38541666666, 1797),
    @(45859.39583333334, 1827),
    @(45859.40625, 1853),
    @(45859.41666666666, 1892),
    @(45859.42708333334, 2023),
    @(45859.4375, 2039),
    @(45859.44791666666, 2057),
    @(45859.45833333334, 2073),
    @(45859.46875, 2135),
    @(45859.47916666666, 2145),
    @(45859.48958333334, 2151),
    @(45859.5, 2154),
    @(45859.51041666666, 2148),
    @(45859.52083333334, 2146),
    @(45859.53125, 2144),
    @(45859.54166666666, 2137),
    @(45859.55208333334, 2078),
    @(45859.5625, 2069),
    @(45859.57291666666, 2055),
    @(45859.58333333334, 2039),
    @(45859.59375, 1931),
    @(45859.60416666666, 1908),
    @(45859.61458333334, 1893),
    @(45859.625, 1871),
    @(45859.63541666666, 1674),
    @(45859.64583333334, 1638),
    @(45859.65625, 1604),
    @(45859.66666666666, 1570),
    @(45859.67708333334, 1318),
    @(45859.6875, 1277),
    @(45859.69791666666, 1242),
    @(45859.70833333334, 1206),
    @(45859.71875, 785),
    @(45859.72916666666, 744),
    @(45859.73958333334, 705),
    @(45859.75, 670),
    @(45859.76041666666, 309),
    @(45859.77083333334, 280),
    @(45859.78125, 253),
    @(45859.79166666666, 233),
    @(45859.80208333334, 58),
    @(45859.8125, 46),
    @(45859.82291666666, 39),
    @(45859.83333333334, 35),
    @(45859.84375, 11),
    @(45859.85416666666, 10),
    @(45859.86458333334, 10),
    @(45859.875, 10),
    @(45859.88541666666, 9),
    @(45859.89583333334, 9),
    @(45859.90625, 8),
    @(45859.91666666666, 8),
    @(45859.92708333334, 1),
    @(45859.9375, 1),
    @(45859.94791666666, 1),
    @(45859.95833333334, 1),
    @(45859.96875, 0),
    @(45859.97916666666, 0),
    @(45859.98958333334, 0),
    @(45860.0, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($row, 2).Value2 = $data[$i][1]
}
